$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A10").Value = "Robert"
$ws.Range("B10").Value = "Stewart"
$ws.Range("C10").Value = "Male"
$ws.Range("D10").Value = $ws.Range("D3").Value2
$ws.Range("D10").NumberFormat = $ws.Range("D3").NumberFormat
$ws.Range("E10").Value = "j.stewart@randatmail.com"
$ws.Range("F10").Value = "319-6123-91"
$ws.Range("G10").Value = "Auditor"
$ws.Range("H10").Value = 98206337
